$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Variables sheet: split the single "position" column into two columns,
# "pivot" (A) and "order" (B). Everything that used to live in columns
# B..O shifts right into C..P.
# ---------------------------------------------------------------------------
$wsVars = $wb.Worksheets.Item("Variables")

# Insert a fresh blank column at A - this pushes the existing "position"
# column (and its column width formatting) from A to B, and B..O to C..P.
$wsVars.Columns.Item(1).Insert()

# Header row
$wsVars.Range("A1").Value = "pivot"
$wsVars.Range("B1").Value = "order"

# Data rows: "h1" -> heading/1, "sN" -> stub/N
$wsVars.Range("A2").Value = "heading"
$wsVars.Range("B2").Value = 1

$wsVars.Range("A3").Value = "stub"
$wsVars.Range("B3").Value = 3
$wsVars.Range("G3").ClearContents()

$wsVars.Range("A4").Value = "stub"
$wsVars.Range("B4").Value = 2
$wsVars.Range("G4").ClearContents()

$wsVars.Range("A5").Value = "stub"
$wsVars.Range("B5").Value = 5
$wsVars.Range("G5").ClearContents()

$wsVars.Range("A6").Value = "stub"
$wsVars.Range("B6").Value = 4
$wsVars.Range("G6").ClearContents()

$wsVars.Range("A7").Value = "stub"
$wsVars.Range("B7").Value = 1
$wsVars.Range("G7").ClearContents()

$wsVars.Range("A8").Value = "stub"
$wsVars.Range("B8").Value = 6
$wsVars.Range("G8").ClearContents()

# Row 9 previously had no "position" value in column A, just a "FIGURES"
# type marker further along the row (now shifted to G9). That marker
# becomes the new pivot value "figures" in column A, and the stray type
# cell is cleared.
$wsVars.Range("A9").Value = "figures"
$wsVars.Range("G9").ClearContents()

# ---------------------------------------------------------------------------
# View/selection bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------------
$wsTable = $wb.Worksheets.Item("Table")
$wsCodelists = $wb.Worksheets.Item("Codelists")

$wsTable.Range("B56").Select()
$wsCodelists.Range("G103").Select()

$wsVars.Activate()
$wsVars.Range("A9").Select()
